# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de) the handback has now completed:
#   - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#     (on the Overview sheet, for every row that showed the old status).
#   - "Latest Target File" (E) / "Latest Handback File" (F) columns get
#     populated with the same file references as "Source File Name" (A) /
#     "Latest Handoff File" (C) respectively, each as a hyperlink.
#   - "Latest Handback DateTime" (G) is stamped with the real handback time.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: flip the status text everywhere it appears ----------
$wb.Worksheets.Item("Overview").Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")

# --- Per-locale detail sheets ---------------------------------------------
function Update-LocaleSheet {
    param($SheetName, $MdBaseUrl, $Row2XlfUrl, $Row3XlfUrl, $Row2Handback, $Row3Handback)

    $ws = $wb.Worksheets.Item($SheetName)

    $aVal2 = $ws.Range("A2").Value2
    $cVal2 = $ws.Range("C2").Value2
    $aVal3 = $ws.Range("A3").Value2
    $cVal3 = $ws.Range("C3").Value2

    # Row 2 ------------------------------------------------------------
    $ws.Range("E2").Value2 = $aVal2
    $ws.Hyperlinks.Add($ws.Range("E2"), $MdBaseUrl + "/e2e/72b5ccf5-3bb2-460a-92a3-d7eb8bba59fb.md", "", "", $aVal2) | Out-Null
    $ws.Range("E2").Style = "HyperLink"

    $ws.Range("F2").Value2 = $cVal2
    $ws.Hyperlinks.Add($ws.Range("F2"), $Row2XlfUrl, "", "", $cVal2) | Out-Null
    $ws.Range("F2").Style = "HyperLink"

    $ws.Range("G2").Value2 = $Row2Handback

    # Row 3 ------------------------------------------------------------
    $ws.Range("E3").Value2 = $aVal3
    $ws.Hyperlinks.Add($ws.Range("E3"), $MdBaseUrl + "/e2e/ccbfa21d-c2ce-41b3-8d42-53855ffb811b.md", "", "", $aVal3) | Out-Null
    $ws.Range("E3").Style = "HyperLink"

    $ws.Range("F3").Value2 = $cVal3
    $ws.Hyperlinks.Add($ws.Range("F3"), $Row3XlfUrl, "", "", $cVal3) | Out-Null
    $ws.Range("F3").Style = "HyperLink"

    $ws.Range("G3").Value2 = $Row3Handback
}

Update-LocaleSheet "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/35dd12414de238a7b57ffd56c8440a3eac9c8e23" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e33354fe9cb4b6fbe2dbb1026d6b1bdcc9122df6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/72b5ccf5-3bb2-460a-92a3-d7eb8bba59fb.d6073054d3a49fec947a960bb8b0e2ad1300d4b1.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e33354fe9cb4b6fbe2dbb1026d6b1bdcc9122df6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ccbfa21d-c2ce-41b3-8d42-53855ffb811b.83adbbd32bb9b7b15e7d17b93b641f434caa6392.zh-cn.xlf" `
    "2016-02-29 04:31:39" `
    "2016-02-29 04:31:39"

Update-LocaleSheet "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/35dd12414de238a7b57ffd56c8440a3eac9c8e23" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18c36f49fb83920d923a808e238912aa76f3d794/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/72b5ccf5-3bb2-460a-92a3-d7eb8bba59fb.d6073054d3a49fec947a960bb8b0e2ad1300d4b1.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/18c36f49fb83920d923a808e238912aa76f3d794/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ccbfa21d-c2ce-41b3-8d42-53855ffb811b.83adbbd32bb9b7b15e7d17b93b641f434caa6392.de-de.xlf" `
    "2016-02-29 04:32:01" `
    "2016-02-29 04:32:01"
